$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new blank rows starting at row 138. This pushes the old
# rows 138-147 down to 147-156 (formulas/relative refs shift automatically).
$ws.Rows("138:146").Insert()

# --- New content for rows 138-140 (three new timing lines, "bam") ---
$ws.Range("B138").Formula = "=C138-B137"
$ws.Range("C138").Value = 207.34
$ws.Range("D138").Value = "bam"

$ws.Range("B139").Formula = "=C139-C138"
$ws.Range("C139").Value = 207.66
$ws.Range("D139").Value = "bam"

$ws.Range("B140").Formula = "=C140-C139"
$ws.Range("C140").Value = 207.97
$ws.Range("D140").Value = "bam"

# --- Row 146 now just holds a gap timing formula between row 140 and 147 ---
$ws.Range("B146").Formula = "=B147-C140"

# Remove the helper cells that Insert() auto-populated with inherited
# formatting/formulas but that have no content in the final layout, so
# those rows/cells don't get serialized as empty (or stale) stubs.
$ws.Range("A138").Clear()
$ws.Range("A139").Clear()
$ws.Range("A140").Clear()
$ws.Range("A141:C145").Clear()
$ws.Range("A146").Clear()
$ws.Range("C146").Clear()
$ws.Range("A147").Clear()

# --- View state: active cell / scroll position ---
$ws.Range("C74").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 133
